$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Préparation" section: add the missing "clip the mixer" step before
#    "Ajouter les ingrédients dans l'ordre." and fix "prends" -> "prend".
#    (Done before touching $d.Tables — doing table work first leaves the
#    Paragraphs collection returning stale text in this runtime.)
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Ajouter les ingrédients dans l’ordre.*") {
        $p.Range.InsertParagraphBefore()
        $newP = $d.Paragraphs.Item($i)
        $newP.Range.Text = "Clipser le mélangeur de la machine."
        break
    }
}

$result = $d.Content.Find.Execute("La cuisson prends environ ", $false, $false, $false, $false, $false, $true, 1, $false, "La cuisson prend environ ", 2)

# ---------------------------------------------------------------------------
# 2) Ingredients table: fix the "7500g" -> "750g" and "1000g" -> "500g" typos
#    in the size-selector header row (row 1).
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

# Row 1, Col 2: "7500g" ("75" + "00g") should read "750g" ("75" + "0g").
$cell2 = $t.Cell(1, 2)
$rng2 = $cell2.Range
$sub2 = $d.Range($rng2.Start + 2, $rng2.Start + 5)
$sub2.Text = "0g"

# Row 1, Col 3: "1000g" should read "500g".
$cell3 = $t.Cell(1, 3)
$rng3 = $cell3.Range
$sub3 = $d.Range($rng3.Start, $rng3.Start + 5)
$sub3.Text = "500g"
